# Update forecast figures after removing Auto ARIMA from the model mix.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# New Prophet/Amazon forecast values for weeks W01..W16 (rows 2-17),
# columns C (Prophet Forecast) through G (Amazon P90 Forecast).
$forecastData = @(
    @(31, 48, 58, 67, 80),
    @(35, 50, 60, 70, 86),
    @(42, 53, 63, 73, 88),
    @(49, 56, 67, 79, 96),
    @(47, 59, 71, 84, 106),
    @(44, 56, 68, 79, 98),
    @(58, 61, 74, 89, 112),
    @(83, 65, 79, 96, 122),
    @(93, 62, 75, 89, 111),
    @(80, 63, 76, 92, 116),
    @(69, 65, 79, 96, 122),
    @(88, 69, 85, 103, 132),
    @(133, 67, 82, 99, 126),
    @(167, 64, 79, 97, 126),
    @(163, 64, 78, 96, 125),
    @(135, 62, 76, 94, 123)
)

$startRow = 2
for ($i = 0; $i -lt $forecastData.Length; $i++) {
    $row = $startRow + $i
    $values = $forecastData[$i]
    $wsForecast.Cells.Item($row, 3).Value = $values[0]
    $wsForecast.Cells.Item($row, 4).Value = $values[1]
    $wsForecast.Cells.Item($row, 5).Value = $values[2]
    $wsForecast.Cells.Item($row, 6).Value = $values[3]
    $wsForecast.Cells.Item($row, 7).Value = $values[4]
}

# --- Sheet 2: "Summary" ----------------------------------------------------
# These "numeric-looking" totals are stored as text on this sheet, so force
# text entry with a leading apostrophe (matches how Excel keeps a typed
# '1317 as text instead of re-interpreting it as a number). Re-apply the
# "Normal" style afterwards so the quote-prefix doesn't leave a stray
# number-format change behind on the cell.
$wsSummary = $wb.Worksheets.Item("Summary")

$summaryUpdates = @{
    "B9"  = "1317"
    "B10" = "389"
    "B11" = "157"
    "B12" = "167"
    "B14" = "31"
}

foreach ($cellRef in $summaryUpdates.Keys) {
    $range = $wsSummary.Range($cellRef)
    $range.Value = "'" + $summaryUpdates[$cellRef]
    $range.Style = "Normal"
}
